# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 12128.333
$ws.Range("I33").Value = 14923.117
$ws.Range("K33").Value = 14923.117
$ws.Range("M33").Value = -14694.117
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 1338.6154
$ws.Range("I40").Value = 1244.4445
$ws.Range("J40").Value = 1550.5
$ws.Range("K40").Value = 1244.4445
$ws.Range("L40").Value = 1550.5
$ws.Range("M40").Value = -1069.4445
$ws.Range("N40").Value = -1900.5
$ws.Range("H42").Value = 324.1111
$ws.Range("I42").Value = 138.93333
$ws.Range("K42").Value = 416.79999
$ws.Range("M42").Value = -186.79999
$ws.Range("H112").Value = 2727
$ws.Range("J112").Value = 3218
$ws.Range("L112").Value = 9654
$ws.Range("N112").Value = -11870
$ws.Range("H116").Value = 29729.646
$ws.Range("I116").Value = 21550
$ws.Range("K116").Value = 21550
$ws.Range("M116").Value = -18108
$ws.Range("H135").Value = 1645.7632
$ws.Range("I135").Value = 1445.8572
$ws.Range("J135").Value = 3978
$ws.Range("K135").Value = 13012.7148
$ws.Range("L135").Value = 35802
$ws.Range("M135").Value = -10477.7148
$ws.Range("N135").Value = -40872
$ws.Range("H137").Value = 23816316
$ws.Range("I137").Value = 33335474
$ws.Range("K137").Value = 100006422
$ws.Range("M137").Value = -100003872

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3707161
$ws.Range("I61").Value = 3470.739
$ws.Range("J61").Value = 25003380
$ws.Range("K61").Value = 3470.739
$ws.Range("L61").Value = 25003380
$ws.Range("M61").Value = -3258.739
$ws.Range("N61").Value = -25003804
$ws.Range("H86").Value = 29000
$ws.Range("J86").Value = 29000
$ws.Range("L86").Value = 29000
$ws.Range("N86").Value = -31372
$ws.Range("H88").Value = 2077.3333
$ws.Range("I88").Value = 921.2
$ws.Range("K88").Value = 921.2
$ws.Range("M88").Value = -515.2
$ws.Range("H89").Value = 29000
$ws.Range("J89").Value = 29000
$ws.Range("L89").Value = 87000
$ws.Range("N89").Value = -98856
$ws.Range("H91").Value = 2077.3333
$ws.Range("I91").Value = 921.2
$ws.Range("K91").Value = 921.2
$ws.Range("M91").Value = 482.8
$ws.Range("H136").Value = 3707161
$ws.Range("I136").Value = 3470.739
$ws.Range("J136").Value = 25003380
$ws.Range("K136").Value = 10412.217
$ws.Range("L136").Value = 75010140
$ws.Range("M136").Value = -7862.217000000001
$ws.Range("N136").Value = -75015240

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 55893.05
$ws.Range("I20").Value = 78469.5
$ws.Range("K20").Value = 78469.5
$ws.Range("M20").Value = -78222.5
$ws.Range("H80").Value = 2101.3333
$ws.Range("J80").Value = 2038.3334
$ws.Range("L80").Value = 2038.3334
$ws.Range("N80").Value = -4034.3334
$ws.Range("H83").Value = 2101.3333
$ws.Range("J83").Value = 2038.3334
$ws.Range("L83").Value = 10191.667
$ws.Range("N83").Value = -20175.667
$ws.Range("H86").Value = 2940.3635
$ws.Range("I86").Value = 2920
$ws.Range("J86").Value = 2994.6667
$ws.Range("K86").Value = 2920
$ws.Range("L86").Value = 2994.6667
$ws.Range("M86").Value = -1797
$ws.Range("N86").Value = -5240.6667
$ws.Range("H89").Value = 2940.3635
$ws.Range("I89").Value = 2920
$ws.Range("J89").Value = 2994.6667
$ws.Range("K89").Value = 14600
$ws.Range("L89").Value = 14973.3335
$ws.Range("M89").Value = -8984
$ws.Range("N89").Value = -26205.3335
$ws.Range("H105").Value = 1923.1875
$ws.Range("I105").Value = 961.1818
$ws.Range("K105").Value = 961.1818
$ws.Range("M105").Value = 785.8182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1869.7646
$ws.Range("I22").Value = 256.57144
$ws.Range("J22").Value = 2999
$ws.Range("K22").Value = 256.57144
$ws.Range("L22").Value = 2999
$ws.Range("M22").Value = 93.42856
$ws.Range("N22").Value = -3699
$ws.Range("H76").Value = 8997.5
$ws.Range("I76").Value = 8997.5
$ws.Range("K76").Value = 8997.5
$ws.Range("M76").Value = -8682.5
$ws.Range("H79").Value = 8997.5
$ws.Range("I79").Value = 8997.5
$ws.Range("K79").Value = 8997.5
$ws.Range("M79").Value = -7905.5
$ws.Range("H105").Value = 8822.5625
$ws.Range("I105").Value = 10683
$ws.Range("J105").Value = 4729.6
$ws.Range("K105").Value = 10683
$ws.Range("L105").Value = 4729.6
$ws.Range("M105").Value = -8936
$ws.Range("N105").Value = -8223.6
$ws.Range("H132").Value = 1403.0741
$ws.Range("I132").Value = 1403.0741
$ws.Range("K132").Value = 4209.2223
$ws.Range("M132").Value = -1679.2223
$ws.Range("H134").Value = 1652.0143
$ws.Range("I134").Value = 1412.2693
$ws.Range("K134").Value = 4236.8079
$ws.Range("M134").Value = -1701.8079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 140600
$ws.Range("J37").Value = 140600
$ws.Range("L37").Value = 421800
$ws.Range("N37").Value = -422024
$ws.Range("H38").Value = 176
$ws.Range("I38").Value = 253.88889
$ws.Range("J38").Value = 125.92857
$ws.Range("K38").Value = 761.6666700000001
$ws.Range("L38").Value = 377.78571
$ws.Range("M38").Value = -414.6666700000001
$ws.Range("N38").Value = -1071.78571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 7000
$ws.Range("I25").Value = 7000
$ws.Range("K25").Value = 7000
$ws.Range("M25").Value = -6471
$ws.Range("H64").Value = 50001
$ws.Range("J64").Value = 50001
$ws.Range("L64").Value = 50001
$ws.Range("N64").Value = -50497
$ws.Range("H67").Value = 50001
$ws.Range("J67").Value = 50001
$ws.Range("L67").Value = 50001
$ws.Range("N67").Value = -51717
$ws.Range("H70").Value = 39709
$ws.Range("I70").Value = 42773.453
$ws.Range("K70").Value = 42773.453
$ws.Range("M70").Value = -42503.453
$ws.Range("H73").Value = 39709
$ws.Range("I73").Value = 42773.453
$ws.Range("K73").Value = 42773.453
$ws.Range("M73").Value = -41837.453
$ws.Range("H80").Value = 2048.1143
$ws.Range("I80").Value = 1914.7241
$ws.Range("K80").Value = 1914.7241
$ws.Range("M80").Value = -916.7240999999999
$ws.Range("H83").Value = 2048.1143
$ws.Range("I83").Value = 1914.7241
$ws.Range("K83").Value = 9573.620499999999
$ws.Range("M83").Value = -4581.620499999999
$ws.Range("H113").Value = 3644.6
$ws.Range("I113").Value = 3305.75
$ws.Range("K113").Value = 3305.75
$ws.Range("M113").Value = -1135.75
$ws.Range("H132").Value = 13223.651
$ws.Range("I132").Value = 13819.756
$ws.Range("J132").Value = 1003.5
$ws.Range("K132").Value = 41459.268
$ws.Range("L132").Value = 3010.5
$ws.Range("M132").Value = -38929.268
$ws.Range("N132").Value = -8070.5
$ws.Range("H135").Value = 75655.73
$ws.Range("J135").Value = 75655.73
$ws.Range("L135").Value = 75655.73
$ws.Range("N135").Value = -85795.73

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2462.8484
$ws.Range("I22").Value = 1103.8
$ws.Range("J22").Value = 3053.739
$ws.Range("K22").Value = 1103.8
$ws.Range("L22").Value = 3053.739
$ws.Range("M22").Value = -808.8
$ws.Range("N22").Value = -3643.739
$ws.Range("H27").Value = 2462.8484
$ws.Range("I27").Value = 1103.8
$ws.Range("J27").Value = 3053.739
$ws.Range("K27").Value = 1103.8
$ws.Range("L27").Value = 3053.739
$ws.Range("M27").Value = -996.8
$ws.Range("N27").Value = -3267.739
$ws.Range("H40").Value = 2853.818
$ws.Range("I40").Value = 2710.3333
$ws.Range("K40").Value = 2710.3333
$ws.Range("M40").Value = -2574.3333
$ws.Range("H61").Value = 5352.643
$ws.Range("I61").Value = 5488.9565
$ws.Range("K61").Value = 5488.9565
$ws.Range("M61").Value = -5286.9565
$ws.Range("H113").Value = 5352.643
$ws.Range("I113").Value = 5488.9565
$ws.Range("K113").Value = 5488.9565
$ws.Range("M113").Value = -3318.9565
$ws.Range("H122").Value = 6580.4546
$ws.Range("I122").Value = 4100
$ws.Range("K122").Value = 12300
$ws.Range("M122").Value = -9850
$ws.Range("H132").Value = 8336707.5
$ws.Range("I132").Value = 9526809
$ws.Range("K132").Value = 28580427
$ws.Range("M132").Value = -28577897

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 58827724
$ws.Range("I81").Value = 4398.125
$ws.Range("K81").Value = 8796.25
$ws.Range("M81").Value = -7735.25
$ws.Range("H84").Value = 58827724
$ws.Range("I84").Value = 4398.125
$ws.Range("K84").Value = 43981.25
$ws.Range("M84").Value = -38677.25
$ws.Range("H122").Value = 45533.074
$ws.Range("I122").Value = 4055.85
$ws.Range("K122").Value = 12167.55
$ws.Range("M122").Value = -9717.549999999999
$ws.Range("H132").Value = 7249158.5
$ws.Range("J132").Value = 3465.2222
$ws.Range("L132").Value = 10395.6666
$ws.Range("N132").Value = -15455.6666
$ws.Range("H136").Value = 6240055
$ws.Range("I136").Value = 2719006.8
$ws.Range("K136").Value = 8157020.399999999
$ws.Range("M136").Value = -8154470.399999999
